$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell values are written in the same order the original authoring session
# entered them (this drives the shared-string table build order), rather
# than strict left-to-right/top-to-bottom order.

# Header row - entered as: Name, email id, username, age, designation, bio,
# work ex, Role, password (columns were rearranged later, hence "age" sits
# at G1 despite being entered before "designation"/"bio"/"work ex").
$ws.Range("A1").Value = "Name"
$ws.Range("B1").Value = "email id"
$ws.Range("C1").Value = "username"
$ws.Range("G1").Value = "age"
$ws.Range("D1").Value = "designation"
$ws.Range("E1").Value = "bio"
$ws.Range("F1").Value = "work ex "
$ws.Range("H1").Value = "Role"
$ws.Range("I1").Value = "password"

# Row 2 core values
$ws.Range("A2").Value = "tb_0"
$ws.Range("B2").Value = "tb@g.com"
$ws.Range("D2").Value = "sdet"
$ws.Range("E2").Value = "coder"
$ws.Range("H2").Value = "Employee"

# Row 3 / Row 4 names
$ws.Range("A3").Value = "tb_1"
$ws.Range("A4").Value = "tb_2"

# skills column added last
$ws.Range("J1").Value = "skills"
$ws.Range("J2").Value = "React.js, Angular"

# Remaining (duplicate-string) cells
$ws.Range("C2").Value = "tb_0"
$ws.Range("C3").Value = "tb_1"
$ws.Range("C4").Value = "tb_2"
$ws.Range("B3").Value = "tb@g.com"
$ws.Range("B4").Value = "tb@g.com"
$ws.Range("D3").Value = "sdet"
$ws.Range("D4").Value = "sdet"
$ws.Range("E3").Value = "coder"
$ws.Range("E4").Value = "coder"
$ws.Range("H3").Value = "Employee"
$ws.Range("H4").Value = "Employee"
$ws.Range("I2").Value = "password"
$ws.Range("I3").Value = "password"
$ws.Range("I4").Value = "password"
$ws.Range("J3").Value = "React.js, Angular"
$ws.Range("J4").Value = "React.js, Angular"

# Numeric cells (no shared-string entries)
$ws.Range("F2").Value = 5
$ws.Range("G2").Value = 23
$ws.Range("F3").Value = 6
$ws.Range("G3").Value = 23
$ws.Range("F4").Value = 7
$ws.Range("G4").Value = 23

# Hyperlinks for the email column
$ws.Hyperlinks.Add($ws.Range("B2"), "mailto:tb@g.com", "", "", "tb@g.com")
$ws.Hyperlinks.Add($ws.Range("B3:B4"), "mailto:tb@g.com", "", "", "tb@g.com")
$ws.Range("B4").Style = "Hyperlink"

# Column widths (values chosen so the runtime's internal pixel-rounding of
# ColumnWidth reproduces the target stored <col width> as closely as possible)
$ws.Columns.Item(1).ColumnWidth = 12.584
$ws.Columns.Item(2).ColumnWidth = 12.917
$ws.Columns.Item(3).ColumnWidth = 17.917
$ws.Columns.Item(4).ColumnWidth = 16.584
$ws.Columns.Item(5).ColumnWidth = 16.917
$ws.Columns.Item(6).ColumnWidth = 17.084
$ws.Columns.Item(7).ColumnWidth = 16.917
$ws.Columns.Item(8).ColumnWidth = 16.751
$ws.Columns.Item(9).ColumnWidth = 16.751
$ws.Columns.Item(10).ColumnWidth = 14.417

$ws.PageSetup.Orientation = 1

$null = $ws.Range("H13").Select()
